# The table on slide 5 (the B1 "types of financial documents" table) had its
# table style switched from the custom "Table_0" style
# ({5140DF42-D9DE-47EE-BC8D-B9D8E97DC5F6}) to the built-in table style
# {932178DA-2DEE-486C-8497-566EF11FD145}.

$p = $ppt.ActivePresentation
$targetStyleId = "{932178DA-2DEE-486C-8497-566EF11FD145}"

$applied = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.StyleId -eq "{5140DF42-D9DE-47EE-BC8D-B9D8E97DC5F6}") {
                $table.ApplyStyle($targetStyleId)
                $applied = $true
            }
        }
    }
}

Write-Output ("Applied=" + $applied)
